# Update "想去人数" (column F) values in the "展览" and "全部类型" sheets
# to reflect the refreshed scrape data (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Row number (in sheet) -> new value for column F
$updates = @{
    2  = 7449
    3  = 7471
    4  = 100
    10 = 137
    13 = 676
    14 = 581
    16 = 35
    19 = 81
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
